# Commit message: "Fruta / hortaliza, semanal"
# This is a weekly data refresh: a new market-report record is inserted as
# row 54 (pushing the existing rows 54-67 down to rows 55-68) on the single
# worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 54; everything below shifts down one row.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new weekly record.
$ws.Range("A54").Value = 7
$ws.Range("B54").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C54").Value = "Ñuble"
$ws.Range("D54").Value = 44617
$ws.Range("E54").Value = 16
$ws.Range("F54").Value = 100112022
$ws.Range("G54").Value = "Arveja Verde"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 80
$ws.Range("K54").Value = 24000
$ws.Range("L54").Value = 25000
$ws.Range("M54").Value = 24500
$ws.Range("N54").Value = "`$/saco 25 kilos"
$ws.Range("O54").Value = "Provincia de Diguillín"
$ws.Range("P54").Value = 980
$ws.Range("Q54").Value = 25
$ws.Range("R54").Value = "Hortaliza"

# Note: Row insertion automatically copies the formatting of the row above,
# so D54 already carries the date number-format style (no need to touch
# .Style explicitly; doing so would strip the inherited number format).
